# Update betting odds / values on Sheet1 to match the 2024-10-17 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8

# Row 3
$ws.Range("V3").Value = 1.67

# Row 4
$ws.Range("V4").Value = 1.67

# Row 5
$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.29

# Row 6
$ws.Range("U6").Value = 1.85
$ws.Range("V6").Value = 1.91

# Row 8
$ws.Range("BC8").Value = 151

# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 2.6
$ws.Range("L10").Value = 4.75
$ws.Range("N10").Value = 8.5
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.73
$ws.Range("X10").Value = 8
$ws.Range("Z10").Value = 15
$ws.Range("AO10").Value = 10
$ws.Range("AQ10").Value = 34

# Row 11
$ws.Range("V11").Value = 1.63

# Row 12
$ws.Range("V12").Value = 1.63

# Row 13
$ws.Range("V13").Value = 1.63
